# "regn, trading updates, crypto updates"
# Adds the 2022-08-08 NAV line to Main, records the RBLX/SGHC closing
# trades (+ their roll-up totals and the new grand total) on Trades, and
# leaves the workbook focused on Trades - matching the author's final
# view state.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Main")
$ws2 = $wb.Worksheets.Item("Trades")

# ---------------------------------------------------------------------
# Main: new NAV row 20 (2022-08-08)
# ---------------------------------------------------------------------

# Copy row 19's formatting down into row 20 first so the new row picks up
# the same number formats / alignment as the rest of the table.
$ws1.Range("B19:G19").Copy()
$ws1.Range("B20:G20").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws1.Range("B20").Value = 44781
$ws1.Range("C20").Value = 45951.47
$ws1.Range("D20").Formula = "=C20-C19"
$ws1.Range("E20").Formula = "=C20-`$C`$4"
$ws1.Range("F20").Formula = "=+C20/C19-1"
$ws1.Range("G20").Formula = "=C20/`$C`$4-1"

# Leave Main's selection parked on the new row (matches the saved view).
$ws1.Range("B20").Select()

# ---------------------------------------------------------------------
# Trades: closed RBLX (-50) and SGHC (+500) positions, new grand total
# ---------------------------------------------------------------------
$ws2.Activate()

# Row 289: RBLX trade
$ws2.Range("B289").Value = "RBLX"
$ws2.Range("C289").Value = "2022-08-08, 10:20:07"
$ws2.Range("D289").Value = -50
$ws2.Range("E289").Value = 49.88
$ws2.Range("F289").Value = 48.9
$ws2.Range("G289").Value = 2494
$ws2.Range("H289").Value = -1.06
$ws2.Range("I289").Value = -2030.5
$ws2.Range("J289").Value = 462.44
$ws2.Range("K289").Value = 49
$ws2.Range("L289").Value = "C;P"

# Row 290: Total RBLX
$ws2.Range("B290").Value = "Total RBLX"
$ws2.Range("D290").Value = -50
$ws2.Range("E290").Value = " "
$ws2.Range("G290").Value = 2494
$ws2.Range("H290").Value = -1.06
$ws2.Range("I290").Value = -2030.5
$ws2.Range("J290").Value = 462.44
$ws2.Range("K290").Value = 49
$ws2.Range("L290").Value = " "
$ws2.Range("L290").NumberFormat = "#,##0.00"

# Row 291: SGHC trade
$ws2.Range("B291").Value = "SGHC"
$ws2.Range("C291").Value = "2022-08-08, 09:36:37"
$ws2.Range("D291").Value = 500
$ws2.Range("E291").Value = 5.3196
$ws2.Range("F291").Value = 5.35
$ws2.Range("G291").Value = -2659.8
$ws2.Range("H291").Value = -2.5
$ws2.Range("I291").Value = 2047.39
$ws2.Range("J291").Value = -614.91
$ws2.Range("K291").Value = 15.2
$ws2.Range("L291").Value = "C;P"

# Row 292: Total SGHC
$ws2.Range("B292").Value = "Total SGHC"
$ws2.Range("D292").Value = 500
$ws2.Range("E292").Value = " "
$ws2.Range("G292").Value = -2659.8
$ws2.Range("H292").Value = -2.5
$ws2.Range("I292").Value = 2047.39
$ws2.Range("J292").Value = -614.91
$ws2.Range("K292").Value = 15.2
$ws2.Range("L292").Value = " "
$ws2.Range("L292").NumberFormat = "#,##0.00"

# Row 293: new grand Total (replaces row 286 as the bottom-most total).
# G293/H293 keep the plain "General" look the earlier total rows have
# (no number-format override), so borrow the already-unstyled look of a
# column-D cell rather than stamping a brand new style index.
$ws2.Range("B293").Value = "Total"
$ws2.Range("G293").Value = -165.8
$ws2.Range("H293").Value = -3.56
$ws2.Range("G293:H293").Style = $ws2.Range("D289").Style
$ws2.Range("I293").Value = 16.89
$ws2.Range("J293").Value = -152.48
$ws2.Range("K293").Value = 64.2
$ws2.Range("L293").Value = " "
$ws2.Range("L293").NumberFormat = "#,##0.00"

# ---------------------------------------------------------------------
# View state: Trades becomes the active / focused sheet, scrolled down
# toward the new rows, Main keeps its selection parked at B20.
# ---------------------------------------------------------------------
$aw = $excel.ActiveWindow
$aw.FreezePanes = $false
$ws2.Range("C3").Select()
$aw.FreezePanes = $true
$ws2.Range("G293:L293").Select()
